$wb = $excel.ActiveWorkbook

# --- Update the "Date" metadata value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-08-02T16:54:26+00:00"

# --- Add two new concept rows on the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Row 16: Level=1, Code=SSOLID, Display=Solid Tumor (somatic)
$concepts.Range("A16").Value = "'1"
$concepts.Range("B16").Value = "SSOLID"
$concepts.Range("C16").Value = "Solid Tumor (somatic)"
$concepts.Range("A15:D15").Copy()
$concepts.Range("A16:D16").PasteSpecial(-4122)

# Row 17: Level=1, Code=SHEMA, Display=Leukemia (somatic)
$concepts.Range("A17").Value = "'1"
$concepts.Range("B17").Value = "SHEMA"
$concepts.Range("C17").Value = "Leukemia (somatic)"
$concepts.Range("A15:D15").Copy()
$concepts.Range("A17:D17").PasteSpecial(-4122)

$excel.CutCopyMode = $false
